$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(21, 8).Value = 10000  # H21: 0 -> 10000
$ws.Cells.Item(21, 9).Value = 10000  # I21: 0 -> 10000
$ws.Cells.Item(21, 11).Value = 10000  # K21: 0 -> 10000
$ws.Cells.Item(21, 13).Value = -9532  # M21: None -> -9532
$ws.Cells.Item(23, 8).Value = 10000  # H23: 0 -> 10000
$ws.Cells.Item(23, 9).Value = 10000  # I23: 0 -> 10000
$ws.Cells.Item(23, 11).Value = 10000  # K23: 0 -> 10000
$ws.Cells.Item(23, 13).Value = -9766  # M23: None -> -9766
$ws.Cells.Item(33, 8).Value = 5985.7144  # H33: 4833.1113 -> 5985.7144
$ws.Cells.Item(33, 9).Value = 150.5  # I33: 140 -> 150.5
$ws.Cells.Item(33, 10).Value = 13766  # J33: 10699.5 -> 13766
$ws.Cells.Item(33, 11).Value = 150.5  # K33: 140 -> 150.5
$ws.Cells.Item(33, 12).Value = 13766  # L33: 10699.5 -> 13766
$ws.Cells.Item(33, 13).Value = 78.5  # M33: 89 -> 78.5
$ws.Cells.Item(33, 14).Value = -14224  # N33: -11157.5 -> -14224
$ws.Cells.Item(34, 8).Value = 1856.7142  # H34: 2200 -> 1856.7142
$ws.Cells.Item(34, 9).Value = 1856.7142  # I34: 2200 -> 1856.7142
$ws.Cells.Item(34, 11).Value = 1856.7142  # K34: 2200 -> 1856.7142
$ws.Cells.Item(34, 13).Value = -1653.7142  # M34: -1997 -> -1653.7142
$ws.Cells.Item(36, 8).Value = 1856.7142  # H36: 2200 -> 1856.7142
$ws.Cells.Item(36, 9).Value = 1856.7142  # I36: 2200 -> 1856.7142
$ws.Cells.Item(36, 11).Value = 1856.7142  # K36: 2200 -> 1856.7142
$ws.Cells.Item(36, 13).Value = -1141.7142  # M36: -1485 -> -1141.7142
$ws.Cells.Item(107, 8).Value = 1521.25  # H107: 1527.9286 -> 1521.25
$ws.Cells.Item(107, 9).Value = 841.94116  # I107: 852.94116 -> 841.94116
$ws.Cells.Item(107, 11).Value = 841.94116  # K107: 852.94116 -> 841.94116
$ws.Cells.Item(107, 13).Value = 1078.05884  # M107: 1067.05884 -> 1078.05884
$ws.Cells.Item(121, 8).Value = 1164  # H121: 746.5 -> 1164
$ws.Cells.Item(121, 10).Value = 1164  # J121: 746.5 -> 1164
$ws.Cells.Item(121, 12).Value = 3492  # L121: 2239.5 -> 3492
$ws.Cells.Item(121, 14).Value = -6986  # N121: -5733.5 -> -6986
$ws.Cells.Item(132, 8).Value = 2784.75  # H132: 2873.8696 -> 2784.75
$ws.Cells.Item(132, 9).Value = 2015.1818  # I132: 2076.1428 -> 2015.1818
$ws.Cells.Item(132, 11).Value = 6045.5454  # K132: 6228.428400000001 -> 6045.5454
$ws.Cells.Item(132, 13).Value = -3515.5454  # M132: -3698.428400000001 -> -3515.5454
$ws.Cells.Item(133, 8).Value = 79998.8  # H133: 79998.28999999999 -> 79998.8
$ws.Cells.Item(133, 10).Value = 79998.8  # J133: 79998.28999999999 -> 79998.8
$ws.Cells.Item(133, 12).Value = 79998.8  # L133: 79998.28999999999 -> 79998.8
$ws.Cells.Item(133, 14).Value = -90118.8  # N133: -90118.28999999999 -> -90118.8
$ws.Cells.Item(138, 8).Value = 4456.0356  # H138: 4714.793 -> 4456.0356
$ws.Cells.Item(138, 10).Value = 4727.6875  # J138: 5153.1177 -> 4727.6875
$ws.Cells.Item(138, 12).Value = 14183.0625  # L138: 15459.3531 -> 14183.0625
$ws.Cells.Item(138, 14).Value = -24463.0625  # N138: -25739.3531 -> -24463.0625

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 1109.1887  # H32: 1109.283 -> 1109.1887
$ws.Cells.Item(32, 9).Value = 1103.6666  # I32: 1103.7646 -> 1103.6666
$ws.Cells.Item(32, 11).Value = 1103.6666  # K32: 1103.7646 -> 1103.6666
$ws.Cells.Item(32, 13).Value = -816.6666  # M32: -816.7646 -> -816.6666
$ws.Cells.Item(97, 8).Value = 2573.2334  # H97: 2569.2 -> 2573.2334
$ws.Cells.Item(97, 9).Value = 2062.375  # I97: 2057.3333 -> 2062.375
$ws.Cells.Item(97, 11).Value = 2062.375  # K97: 2057.3333 -> 2062.375
$ws.Cells.Item(97, 13).Value = -1566.375  # M97: -1561.3333 -> -1566.375
$ws.Cells.Item(110, 8).Value = 5776.25  # H110: 2903 -> 5776.25
$ws.Cells.Item(110, 9).Value = 3874.75  # I110: 1797.1875 -> 3874.75
$ws.Cells.Item(110, 10).Value = 7677.75  # J110: 6441.6 -> 7677.75
$ws.Cells.Item(110, 11).Value = 3874.75  # K110: 1797.1875 -> 3874.75
$ws.Cells.Item(110, 12).Value = 7677.75  # L110: 6441.6 -> 7677.75
$ws.Cells.Item(110, 13).Value = -1829.75  # M110: 247.8125 -> -1829.75
$ws.Cells.Item(110, 14).Value = -11767.75  # N110: -10531.6 -> -11767.75
$ws.Cells.Item(122, 8).Value = 3128.4866  # H122: 3165.838 -> 3128.4866
$ws.Cells.Item(122, 9).Value = 2771.24  # I122: 2826.44 -> 2771.24
$ws.Cells.Item(122, 10).Value = 3872.75  # J122: 3872.9167 -> 3872.75
$ws.Cells.Item(122, 11).Value = 8313.719999999999  # K122: 8479.32 -> 8313.719999999999
$ws.Cells.Item(122, 12).Value = 11618.25  # L122: 11618.7501 -> 11618.25
$ws.Cells.Item(122, 13).Value = -5863.719999999999  # M122: -6029.32 -> -5863.719999999999
$ws.Cells.Item(122, 14).Value = -16518.25  # N122: -16518.7501 -> -16518.25
$ws.Cells.Item(132, 8).Value = 10728.467  # H132: 9470.789000000001 -> 10728.467
$ws.Cells.Item(132, 9).Value = 7578.8  # I132: 6323.5557 -> 7578.8
$ws.Cells.Item(132, 11).Value = 22736.4  # K132: 18970.6671 -> 22736.4
$ws.Cells.Item(132, 13).Value = -20206.4  # M132: -16440.6671 -> -20206.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(80, 8).Value = 502.5  # H80: 504.17648 -> 502.5
$ws.Cells.Item(80, 9).Value = 451.5  # I80: 401.14285 -> 451.5
$ws.Cells.Item(80, 10).Value = 533.1  # J80: 576.3 -> 533.1
$ws.Cells.Item(80, 11).Value = 451.5  # K80: 401.14285 -> 451.5
$ws.Cells.Item(80, 12).Value = 533.1  # L80: 576.3 -> 533.1
$ws.Cells.Item(80, 13).Value = 546.5  # M80: 596.85715 -> 546.5
$ws.Cells.Item(80, 14).Value = -2529.1  # N80: -2572.3 -> -2529.1
$ws.Cells.Item(83, 8).Value = 502.5  # H83: 504.17648 -> 502.5
$ws.Cells.Item(83, 9).Value = 451.5  # I83: 401.14285 -> 451.5
$ws.Cells.Item(83, 10).Value = 533.1  # J83: 576.3 -> 533.1
$ws.Cells.Item(83, 11).Value = 2257.5  # K83: 2005.71425 -> 2257.5
$ws.Cells.Item(83, 12).Value = 2665.5  # L83: 2881.5 -> 2665.5
$ws.Cells.Item(83, 13).Value = 2734.5  # M83: 2986.28575 -> 2734.5
$ws.Cells.Item(83, 14).Value = -12649.5  # N83: -12865.5 -> -12649.5
$ws.Cells.Item(86, 8).Value = 4495.6924  # H86: 4518.769 -> 4495.6924
$ws.Cells.Item(86, 9).Value = 2343.8  # I86: 2373.8 -> 2343.8
$ws.Cells.Item(86, 11).Value = 2343.8  # K86: 2373.8 -> 2343.8
$ws.Cells.Item(86, 13).Value = -1220.8  # M86: -1250.8 -> -1220.8
$ws.Cells.Item(89, 8).Value = 4495.6924  # H89: 4518.769 -> 4495.6924
$ws.Cells.Item(89, 9).Value = 2343.8  # I89: 2373.8 -> 2343.8
$ws.Cells.Item(89, 11).Value = 11719  # K89: 11869 -> 11719
$ws.Cells.Item(89, 13).Value = -6103  # M89: -6253 -> -6103
$ws.Cells.Item(94, 8).Value = 903.25  # H94: 802.4 -> 903.25
$ws.Cells.Item(94, 9).Value = 903.25  # I94: 802.4 -> 903.25
$ws.Cells.Item(94, 11).Value = 903.25  # K94: 802.4 -> 903.25
$ws.Cells.Item(94, 13).Value = -452.25  # M94: -351.4 -> -452.25
$ws.Cells.Item(99, 8).Value = 3095.5417  # H99: 3138.3914 -> 3095.5417
$ws.Cells.Item(99, 9).Value = 2892.65  # I99: 2933.842 -> 2892.65
$ws.Cells.Item(99, 11).Value = 2892.65  # K99: 2933.842 -> 2892.65
$ws.Cells.Item(99, 13).Value = -1394.65  # M99: -1435.842 -> -1394.65

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(2, 8).Value = 2620  # H2: 2619.8 -> 2620
$ws.Cells.Item(2, 9).Value = 1500  # I2: 1499.5 -> 1500
$ws.Cells.Item(2, 11).Value = 1500  # K2: 1499.5 -> 1500
$ws.Cells.Item(2, 13).Value = -1387  # M2: -1386.5 -> -1387
$ws.Cells.Item(31, 8).Value = 28106.4  # H31: 29961.477 -> 28106.4
$ws.Cells.Item(31, 9).Value = 3142.5833  # I31: 3430.4 -> 3142.5833
$ws.Cells.Item(31, 10).Value = 37184.152  # J31: 38252.438 -> 37184.152
$ws.Cells.Item(31, 11).Value = 3142.5833  # K31: 3430.4 -> 3142.5833
$ws.Cells.Item(31, 12).Value = 37184.152  # L31: 38252.438 -> 37184.152
$ws.Cells.Item(31, 13).Value = -2847.5833  # M31: -3135.4 -> -2847.5833
$ws.Cells.Item(31, 14).Value = -37774.152  # N31: -38842.438 -> -37774.152
$ws.Cells.Item(34, 8).Value = 28106.4  # H34: 29961.477 -> 28106.4
$ws.Cells.Item(34, 9).Value = 3142.5833  # I34: 3430.4 -> 3142.5833
$ws.Cells.Item(34, 10).Value = 37184.152  # J34: 38252.438 -> 37184.152
$ws.Cells.Item(34, 11).Value = 3142.5833  # K34: 3430.4 -> 3142.5833
$ws.Cells.Item(34, 12).Value = 37184.152  # L34: 38252.438 -> 37184.152
$ws.Cells.Item(34, 13).Value = -2940.5833  # M34: -3228.4 -> -2940.5833
$ws.Cells.Item(34, 14).Value = -37588.152  # N34: -38656.438 -> -37588.152
$ws.Cells.Item(59, 8).Value = 35000  # H59: 50000 -> 35000
$ws.Cells.Item(59, 9).Value = 0  # I59: 50000 -> 0
$ws.Cells.Item(59, 10).Value = 35000  # J59: 0 -> 35000
$ws.Cells.Item(59, 11).Value = 0  # K59: 50000 -> 0
$ws.Cells.Item(59, 12).Value = 35000  # L59: 0 -> 35000
$ws.Cells.Item(59, 13).ClearContents()  # M59: -48855 -> (cleared)
$ws.Cells.Item(59, 14).Value = -37290  # N59: None -> -37290
$ws.Cells.Item(132, 8).Value = 4508.4546  # H132: 4544.591 -> 4508.4546
$ws.Cells.Item(132, 9).Value = 3222.7693  # I132: 3283.923 -> 3222.7693
$ws.Cells.Item(132, 11).Value = 9668.3079  # K132: 9851.769 -> 9668.3079
$ws.Cells.Item(132, 13).Value = -7138.3079  # M132: -7321.769 -> -7138.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 1313.35  # H97: 1238.35 -> 1313.35
$ws.Cells.Item(97, 9).Value = 962.5714  # I97: 855.4286 -> 962.5714
$ws.Cells.Item(97, 11).Value = 962.5714  # K97: 855.4286 -> 962.5714
$ws.Cells.Item(97, 13).Value = -466.5714  # M97: -359.4286 -> -466.5714
$ws.Cells.Item(122, 8).Value = 8365.931  # H122: 7743.2812 -> 8365.931
$ws.Cells.Item(122, 9).Value = 7523.52  # I122: 6902.1787 -> 7523.52
$ws.Cells.Item(122, 11).Value = 22570.56  # K122: 20706.5361 -> 22570.56
$ws.Cells.Item(122, 13).Value = -20120.56  # M122: -18256.5361 -> -20120.56
$ws.Cells.Item(132, 8).Value = 6622.4707  # H132: 4844.1665 -> 6622.4707
$ws.Cells.Item(132, 9).Value = 4202.5  # I132: 3414.9583 -> 4202.5
$ws.Cells.Item(132, 10).Value = 8773.556  # J132: 7702.5835 -> 8773.556
$ws.Cells.Item(132, 11).Value = 12607.5  # K132: 10244.8749 -> 12607.5
$ws.Cells.Item(132, 12).Value = 26320.668  # L132: 23107.7505 -> 26320.668
$ws.Cells.Item(132, 13).Value = -10077.5  # M132: -7714.874899999999 -> -10077.5
$ws.Cells.Item(132, 14).Value = -31380.668  # N132: -28167.7505 -> -31380.668

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value = 33312.383  # H22: 31808.818 -> 33312.383
$ws.Cells.Item(22, 9).Value = 50175.375  # I22: 44626.332 -> 50175.375
$ws.Cells.Item(22, 11).Value = 50175.375  # K22: 44626.332 -> 50175.375
$ws.Cells.Item(22, 13).Value = -49880.375  # M22: -44331.332 -> -49880.375
$ws.Cells.Item(27, 8).Value = 33312.383  # H27: 31808.818 -> 33312.383
$ws.Cells.Item(27, 9).Value = 50175.375  # I27: 44626.332 -> 50175.375
$ws.Cells.Item(27, 11).Value = 50175.375  # K27: 44626.332 -> 50175.375
$ws.Cells.Item(27, 13).Value = -50068.375  # M27: -44519.332 -> -50068.375
$ws.Cells.Item(46, 8).Value = 3625.1538  # H46: 3824.9167 -> 3625.1538
$ws.Cells.Item(46, 9).Value = 1403.7142  # I46: 1433 -> 1403.7142
$ws.Cells.Item(46, 11).Value = 1403.7142  # K46: 1433 -> 1403.7142
$ws.Cells.Item(46, 13).Value = -1215.7142  # M46: -1245 -> -1215.7142
$ws.Cells.Item(82, 8).Value = 7375.25  # H82: 11250.167 -> 7375.25
$ws.Cells.Item(82, 10).Value = 0  # J82: 19000 -> 0
$ws.Cells.Item(82, 12).Value = 0  # L82: 19000 -> 0
$ws.Cells.Item(82, 14).ClearContents()  # N82: -19722 -> (cleared)
$ws.Cells.Item(85, 8).Value = 7375.25  # H85: 11250.167 -> 7375.25
$ws.Cells.Item(85, 10).Value = 0  # J85: 19000 -> 0
$ws.Cells.Item(85, 12).Value = 0  # L85: 19000 -> 0
$ws.Cells.Item(85, 14).ClearContents()  # N85: -21496 -> (cleared)
$ws.Cells.Item(93, 8).Value = 12674.913  # H93: 12180.125 -> 12674.913
$ws.Cells.Item(93, 9).Value = 12125.267  # I93: 11417.4375 -> 12125.267
$ws.Cells.Item(93, 11).Value = 12125.267  # K93: 11417.4375 -> 12125.267
$ws.Cells.Item(93, 13).Value = -10877.267  # M93: -10169.4375 -> -10877.267
$ws.Cells.Item(100, 8).Value = 3987.25  # H100: 4385.5713 -> 3987.25
$ws.Cells.Item(100, 9).Value = 1842  # I100: 1949.1666 -> 1842
$ws.Cells.Item(100, 11).Value = 1842  # K100: 1949.1666 -> 1842
$ws.Cells.Item(100, 13).Value = -1301  # M100: -1408.1666 -> -1301
$ws.Cells.Item(122, 8).Value = 8613  # H122: 8989.666999999999 -> 8613
$ws.Cells.Item(122, 9).Value = 4975.25  # I122: 4978 -> 4975.25
$ws.Cells.Item(122, 10).Value = 12250.75  # J122: 13001.333 -> 12250.75
$ws.Cells.Item(122, 11).Value = 14925.75  # K122: 14934 -> 14925.75
$ws.Cells.Item(122, 12).Value = 36752.25  # L122: 39003.999 -> 36752.25
$ws.Cells.Item(122, 13).Value = -12475.75  # M122: -12484 -> -12475.75
$ws.Cells.Item(122, 14).Value = -41652.25  # N122: -43903.999 -> -41652.25
$ws.Cells.Item(132, 8).Value = 6929.0713  # H132: 4600.7144 -> 6929.0713
$ws.Cells.Item(132, 9).Value = 4285.7144  # I132: 3290.7334 -> 4285.7144
$ws.Cells.Item(132, 10).Value = 9572.429  # J132: 6112.231 -> 9572.429
$ws.Cells.Item(132, 11).Value = 12857.1432  # K132: 9872.200199999999 -> 12857.1432
$ws.Cells.Item(132, 12).Value = 28717.287  # L132: 18336.693 -> 28717.287
$ws.Cells.Item(132, 13).Value = -10327.1432  # M132: -7342.200199999999 -> -10327.1432
$ws.Cells.Item(132, 14).Value = -33777.287  # N132: -23396.693 -> -33777.287

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(6, 8).Value = 2500  # H6: 0 -> 2500
$ws.Cells.Item(6, 10).Value = 2500  # J6: 0 -> 2500
$ws.Cells.Item(6, 12).Value = 2500  # L6: 0 -> 2500
$ws.Cells.Item(6, 14).Value = -2730  # N6: None -> -2730
$ws.Cells.Item(13, 8).Value = 8913.571  # H13: 9100 -> 8913.571
$ws.Cells.Item(13, 10).Value = 8879  # J13: 9200 -> 8879
$ws.Cells.Item(13, 12).Value = 8879  # L13: 9200 -> 8879
$ws.Cells.Item(13, 14).Value = -9159  # N13: -9480 -> -9159
$ws.Cells.Item(34, 8).Value = 0  # H34: 11024.667 -> 0
$ws.Cells.Item(34, 9).Value = 0  # I34: 11024 -> 0
$ws.Cells.Item(34, 10).Value = 0  # J34: 11025 -> 0
$ws.Cells.Item(34, 11).Value = 0  # K34: 11024 -> 0
$ws.Cells.Item(34, 12).Value = 0  # L34: 11025 -> 0
$ws.Cells.Item(34, 13).ClearContents()  # M34: -10821 -> (cleared)
$ws.Cells.Item(34, 14).ClearContents()  # N34: -11431 -> (cleared)
$ws.Cells.Item(81, 8).Value = 4786.9  # H81: 5109.3335 -> 4786.9
$ws.Cells.Item(81, 9).Value = 3207.4443  # I81: 3372.75 -> 3207.4443
$ws.Cells.Item(81, 11).Value = 6414.8886  # K81: 6745.5 -> 6414.8886
$ws.Cells.Item(81, 13).Value = -5353.8886  # M81: -5684.5 -> -5353.8886
$ws.Cells.Item(84, 8).Value = 4786.9  # H84: 5109.3335 -> 4786.9
$ws.Cells.Item(84, 9).Value = 3207.4443  # I84: 3372.75 -> 3207.4443
$ws.Cells.Item(84, 11).Value = 32074.443  # K84: 33727.5 -> 32074.443
$ws.Cells.Item(84, 13).Value = -26770.443  # M84: -28423.5 -> -26770.443
$ws.Cells.Item(122, 8).Value = 2402.5715  # H122: 2545.7368 -> 2402.5715
$ws.Cells.Item(122, 9).Value = 1565.5555  # I122: 1630.9375 -> 1565.5555
$ws.Cells.Item(122, 11).Value = 4696.666499999999  # K122: 4892.8125 -> 4696.666499999999
$ws.Cells.Item(122, 13).Value = -2246.666499999999  # M122: -2442.8125 -> -2246.666499999999
$ws.Cells.Item(132, 8).Value = 11824.5  # H132: 9001.913 -> 11824.5
$ws.Cells.Item(132, 9).Value = 5198.909  # I132: 4168.8887 -> 5198.909
$ws.Cells.Item(132, 11).Value = 15596.727  # K132: 12506.6661 -> 15596.727
$ws.Cells.Item(132, 13).Value = -13066.727  # M132: -9976.666100000002 -> -13066.727
$ws.Cells.Item(136, 8).Value = 3193.5918  # H136: 3309.9348 -> 3193.5918
$ws.Cells.Item(136, 9).Value = 2857.1904  # I136: 2950.325 -> 2857.1904
$ws.Cells.Item(136, 10).Value = 5212  # J136: 5707.3335 -> 5212
$ws.Cells.Item(136, 11).Value = 8571.5712  # K136: 8850.974999999999 -> 8571.5712
$ws.Cells.Item(136, 12).Value = 15636  # L136: 17122.0005 -> 15636
$ws.Cells.Item(136, 13).Value = -6021.5712  # M136: -6300.974999999999 -> -6021.5712
$ws.Cells.Item(136, 14).Value = -20736  # N136: -22222.0005 -> -20736
